$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author prepended five earlier historical data points (1523-1651) to the
# breeding-probability table and refreshed the probability figures for every
# existing year. Make room for the new rows first, then write the full,
# renumbered data block (years 1523-2012, rows 2-29).
$ws.Range("A2:B6").EntireRow.Insert()

$years = @(1523,1542,1561,1601,1651,1701,1802,1812,1822,1832,1842,1852,1862,1872,1882,1892,1902,1912,1922,1932,1942,1952,1962,1972,1982,1992,2002,2012)
$probs = @(0.35,0.35,0.31,0.29,0.26,0.25,0.28,0.28,0.26,0.25,0.22,0.25,0.22,0.21,0.19,0.18,0.13,0.14,0.12,0.07,0.08,0.09,0.1,0.07,0.1,0.09,0.09,0.09)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value2 = $years[$i]
    $ws.Cells.Item($row, 2).Value2 = $probs[$i]
}

# Match the author's final selection (cell B30, just past the new data block).
[void]$ws.Range("B30").Select()

# Keep the scatter chart's X axis and position in sync with the now-larger
# data range: the earliest year moved from 1701 to 1523, and the chart itself
# shifted down by five rows (75pt at the sheet's 15pt row height) so it still
# sits just below the data table.
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$xAxis = $chart.Axes(1)
$xAxis.MinimumScale = 1523
$co.Top = $co.Top + 75
